$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2249.6667
$ws.Range("I19").Value = 1831.6666
$ws.Range("J19").Value = 2458.6667
$ws.Range("K19").Value = 1831.6666
$ws.Range("L19").Value = 2458.6667
$ws.Range("M19").Value = -1656.6666
$ws.Range("N19").Value = -2808.6667
$ws.Range("H40").Value = 3809.0312
$ws.Range("I40").Value = 3361.875
$ws.Range("J40").Value = 4256.1875
$ws.Range("K40").Value = 3361.875
$ws.Range("L40").Value = 4256.1875
$ws.Range("M40").Value = -3186.875
$ws.Range("N40").Value = -4606.1875
$ws.Range("H43").Value = 2562.1785
$ws.Range("I43").Value = 1344.8096
$ws.Range("K43").Value = 1344.8096
$ws.Range("M43").Value = -1275.8096
$ws.Range("H51").Value = 37945.453
$ws.Range("I51").Value = 150000.0
$ws.Range("K51").Value = 150000.0
$ws.Range("M51").Value = -149516.0
$ws.Range("H55").Value = 1856.6666
$ws.Range("I55").Value = 175.88889
$ws.Range("J55").Value = 6899.0
$ws.Range("K55").Value = 175.88889
$ws.Range("L55").Value = 6899.0
$ws.Range("M55").Value = 38.11111
$ws.Range("N55").Value = -7327.0
$ws.Range("H98").Value = 37038520.0
$ws.Range("I98").Value = 37038520.0
$ws.Range("K98").Value = 37038520.0
$ws.Range("M98").Value = -37037022.0
$ws.Range("H116").Value = 5737.0713
$ws.Range("I116").Value = 5785.5713
$ws.Range("J116").Value = 5688.5713
$ws.Range("K116").Value = 5785.5713
$ws.Range("L116").Value = 5688.5713
$ws.Range("M116").Value = -2343.5713
$ws.Range("N116").Value = -12572.5713
$ws.Range("H122").Value = 37038520.0
$ws.Range("I122").Value = 37038520.0
$ws.Range("K122").Value = 111115560.0
$ws.Range("M122").Value = -111113110.0

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38545576.0
$ws.Range("I61").Value = 83334860.0
$ws.Range("K61").Value = 83334860.0
$ws.Range("M61").Value = -83334648.0
$ws.Range("H110").Value = 5732.8335
$ws.Range("I110").Value = 2879.4
$ws.Range("J110").Value = 20000.0
$ws.Range("K110").Value = 2879.4
$ws.Range("L110").Value = 20000.0
$ws.Range("M110").Value = -834.4000000000001
$ws.Range("N110").Value = -24090.0
$ws.Range("H114").Value = 0.0
$ws.Range("J114").Value = 0.0
$ws.Range("L114").Value = 0.0
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 1444.625
$ws.Range("I122").Value = 1093.0
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 3279.0
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -829.0
$ws.Range("N122").Value = -12398.5
$ws.Range("H132").Value = 4804.2354
$ws.Range("I132").Value = 2465.32
$ws.Range("K132").Value = 7395.960000000001
$ws.Range("M132").Value = -4865.960000000001
$ws.Range("H136").Value = 38545576.0
$ws.Range("I136").Value = 83334860.0
$ws.Range("K136").Value = 250004580.0
$ws.Range("M136").Value = -250002030.0

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9347.0
$ws.Range("I20").Value = 9347.0
$ws.Range("K20").Value = 9347.0
$ws.Range("M20").Value = -9100.0
$ws.Range("H86").Value = 1997.9642
$ws.Range("I86").Value = 1716.5625
$ws.Range("J86").Value = 2373.1667
$ws.Range("K86").Value = 1716.5625
$ws.Range("L86").Value = 2373.1667
$ws.Range("M86").Value = -593.5625
$ws.Range("N86").Value = -4619.1667
$ws.Range("H89").Value = 1997.9642
$ws.Range("I89").Value = 1716.5625
$ws.Range("J89").Value = 2373.1667
$ws.Range("K89").Value = 8582.8125
$ws.Range("L89").Value = 11865.8335
$ws.Range("M89").Value = -2966.8125
$ws.Range("N89").Value = -23097.8335
$ws.Range("H105").Value = 2343.0
$ws.Range("I105").Value = 2140.8333
$ws.Range("J105").Value = 2949.5
$ws.Range("K105").Value = 2140.8333
$ws.Range("L105").Value = 2949.5
$ws.Range("M105").Value = -393.8332999999998
$ws.Range("N105").Value = -6443.5
$ws.Range("H107").Value = 2018.7333
$ws.Range("I107").Value = 1773.25
$ws.Range("J107").Value = 3000.6667
$ws.Range("K107").Value = 1773.25
$ws.Range("L107").Value = 3000.6667
$ws.Range("M107").Value = 146.75
$ws.Range("N107").Value = -6840.6667
$ws.Range("H134").Value = 36026.1
$ws.Range("I134").Value = 1409.75
$ws.Range("J134").Value = 202184.6
$ws.Range("K134").Value = 4229.25
$ws.Range("L134").Value = 606553.8
$ws.Range("M134").Value = -1694.25
$ws.Range("N134").Value = -611623.8

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 75357.164
$ws.Range("J108").Value = 75357.164
$ws.Range("L108").Value = 75357.164
$ws.Range("N108").Value = -83037.164
$ws.Range("H132").Value = 3098.0476
$ws.Range("I132").Value = 3187.2222
$ws.Range("K132").Value = 9561.6666
$ws.Range("M132").Value = -7031.6666

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 697.35297
$ws.Range("I5").Value = 597.3333
$ws.Range("K5").Value = 1791.9999
$ws.Range("M5").Value = -1679.9999
$ws.Range("H11").Value = 469.89474
$ws.Range("I11").Value = 345.5
$ws.Range("J11").Value = 1133.3334
$ws.Range("K11").Value = 1036.5
$ws.Range("L11").Value = 3400.0002
$ws.Range("M11").Value = -896.5
$ws.Range("N11").Value = -3680.0002
$ws.Range("H57").Value = 0.0
$ws.Range("J57").Value = 0.0
$ws.Range("L57").Value = 0.0
$ws.Range("N57").ClearContents()
$ws.Range("H80").Value = 3383.3
$ws.Range("J80").Value = 3432.3333
$ws.Range("L80").Value = 10296.9999
$ws.Range("N80").Value = -12168.9999
$ws.Range("H83").Value = 3383.3
$ws.Range("J83").Value = 3432.3333
$ws.Range("L83").Value = 30890.9997
$ws.Range("N83").Value = -40250.9997
$ws.Range("H135").Value = 697.35297
$ws.Range("I135").Value = 597.3333
$ws.Range("K135").Value = 5375.9997
$ws.Range("M135").Value = -2840.9997

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2057.2917
$ws.Range("I102").Value = 2084.5715
$ws.Range("K102").Value = 2084.5715
$ws.Range("M102").Value = -462.5715
$ws.Range("H126").Value = 2278.125
$ws.Range("I126").Value = 2146.4
$ws.Range("J126").Value = 2497.6667
$ws.Range("K126").Value = 6439.200000000001
$ws.Range("L126").Value = 7493.000100000001
$ws.Range("M126").Value = -3969.200000000001
$ws.Range("N126").Value = -12433.0001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 19835.334
$ws.Range("I40").Value = 2915.3333
$ws.Range("J40").Value = 22655.334
$ws.Range("K40").Value = 2915.3333
$ws.Range("L40").Value = 22655.334
$ws.Range("M40").Value = -2779.3333
$ws.Range("N40").Value = -22927.334
$ws.Range("H55").Value = 52631896.0
$ws.Range("I55").Value = 58823860.0
$ws.Range("J55").Value = 201.5
$ws.Range("K55").Value = 58823860.0
$ws.Range("L55").Value = 201.5
$ws.Range("M55").Value = -58823687.0
$ws.Range("N55").Value = -547.5
$ws.Range("H68").Value = 3896.5715
$ws.Range("I68").Value = 3406.2
$ws.Range("K68").Value = 3406.2
$ws.Range("M68").Value = -2657.2
$ws.Range("H71").Value = 3896.5715
$ws.Range("I71").Value = 3406.2
$ws.Range("K71").Value = 17031.0
$ws.Range("M71").Value = -13287.0
$ws.Range("H132").Value = 54987.867
$ws.Range("J132").Value = 251949.75
$ws.Range("L132").Value = 755849.25
$ws.Range("N132").Value = -760909.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 15050699.0
$ws.Range("J75").Value = 15050699.0
$ws.Range("L75").Value = 15050699.0
$ws.Range("N75").Value = -15052571.0
$ws.Range("H78").Value = 15050699.0
$ws.Range("J78").Value = 15050699.0
$ws.Range("L78").Value = 45152097.0
$ws.Range("N78").Value = -45161457.0
$ws.Range("H81").Value = 2995.0
$ws.Range("I81").Value = 2990.0
$ws.Range("K81").Value = 5980.0
$ws.Range("M81").Value = -4919.0
$ws.Range("H84").Value = 2995.0
$ws.Range("I84").Value = 2990.0
$ws.Range("K84").Value = 29900.0
$ws.Range("M84").Value = -24596.0
$ws.Range("H96").Value = 3568.6
$ws.Range("I96").Value = 1963.0
$ws.Range("K96").Value = 1963.0
$ws.Range("M96").Value = -590.0
$ws.Range("H116").Value = 128000.0
$ws.Range("J116").Value = 128000.0
$ws.Range("L116").Value = 128000.0
$ws.Range("N116").Value = -137178.0
$ws.Range("H132").Value = 2476.0
$ws.Range("I132").Value = 1841.0834
$ws.Range("K132").Value = 5523.2502
$ws.Range("M132").Value = -2993.2502
